$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.583.45"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "2.631.04"
$ws.Range("E3").Value = "  -1.80%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.34"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.74"
$ws.Range("E6").Value = "  +0.95%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -2.21%  "
$ws.Range("D9").Value = "2.629.72"
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.23"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.71"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").Value = "3.106.89"
$ws.Range("E15").Value = "  -1.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000183"
$ws.Range("E16").Value = "  -1.25%  "
$ws.Range("D17").Value = "67.430.00"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "2.638.32"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.04"
$ws.Range("E19").Value = "  +2.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.06"
$ws.Range("E20").Value = "  +4.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "357.64"
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("E22").Value = "  -1.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.68"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.95"
$ws.Range("E24").Value = "  -4.04%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.33"
$ws.Range("E26").Value = "  +3.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "69.68"
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("D28").Value = "2.758.98"
$ws.Range("E28").Value = "  +1.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "548.18"
$ws.Range("E31").Value = "  -1.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.95"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.36"
$ws.Range("E33").Value = "  -2.76%  "
$ws.Range("E35").Value = "  +4.34%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.51"
$ws.Range("E37").Value = "  -3.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.85"
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.03"
$ws.Range("E39").Value = "  -2.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.367"
$ws.Range("E40").Value = "  -1.79%  "
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.30"
$ws.Range("E42").Value = "  +1.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.24"
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.43"
$ws.Range("E45").Value = "  -3.76%  "
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "152.93"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.580"
$ws.Range("E48").Value = "  -1.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.80"
$ws.Range("E49").Value = "  -1.31%  "
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("E51").Value = "  -1.21%  "
